# Re-pulled data: push updated "dSF" (column F) values for the
# stat-table rows. Values were refreshed from source and therefore
# differ from the previously stored (often zero) placeholders.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$dsf = @{
    2  = 1
    3  = -2
    4  = 1
    5  = -1
    6  = -1
    7  = -2
    8  = 1
    9  = 5
    10 = 5
    11 = 3
    12 = -3
    13 = 6
    14 = 2
    15 = 7
    16 = 2
    17 = -4
    18 = 1
    19 = -3
    20 = -2
    21 = 2
    23 = -1
    24 = -2
    26 = 1
}

foreach ($row in $dsf.Keys) {
    $ws.Cells.Item($row, 6).Value = $dsf[$row]
}
